# "Pushed additional Audio Files"
#
# - Append two new hyperlink paragraphs (matching the style of the
#   existing Freesound links: a hyperlink run followed by a trailing
#   space run) at the end of the "Source List" document.
# - The "_GoBack" bookmark (left over from the last edit position of the
#   previous save) moves from the heading paragraph to the new, final
#   empty paragraph - this is normal Word behaviour: it always marks the
#   location of the most recent edit.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Drop the stale "_GoBack" bookmark that currently sits on the
#    "Source List" heading paragraph. (_GoBack is a hidden bookmark, so
#    it is reached by name rather than through the Bookmarks collection
#    enumerator.)
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Helper pattern used twice below: turn the (currently last, empty)
# paragraph into "<hyperlink> " and then split off a fresh empty
# paragraph after it via Find/Replace with "^p" (this keeps the new
# paragraph free of any stray empty run, unlike Range.InsertParagraphAfter).
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 2. New link #1: anderz000
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$d.Hyperlinks.Add($r, "https://freesound.org/people/anderz000/sounds/204310/", $null, $null, `
    "https://freesound.org/people/anderz000/sounds/204310/")

$p = $d.Paragraphs.Last
$pr = $p.Range
$pr.Collapse(0)
$pr.Font.Reset()
$pr.InsertAfter(" ")

$p = $d.Paragraphs.Last
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter("@@SplitMarker@@")
$d.Content.Find.Execute("@@SplitMarker@@", $true, $false, $false, $false, $false, `
    $true, 1, $false, "^p", 2)

# ---------------------------------------------------------------------
# 3. New link #2: LittleRobotSoundFactory
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$d.Hyperlinks.Add($r, "https://freesound.org/people/LittleRobotSoundFactory/sounds/270404/", $null, $null, `
    "https://freesound.org/people/LittleRobotSoundFactory/sounds/270404/")

$p = $d.Paragraphs.Last
$pr = $p.Range
$pr.Collapse(0)
$pr.Font.Reset()
$pr.InsertAfter(" ")

$p = $d.Paragraphs.Last
$pr = $p.Range
$pr.Collapse(0)
$pr.InsertAfter("@@SplitMarker@@")
$d.Content.Find.Execute("@@SplitMarker@@", $true, $false, $false, $false, $false, `
    $true, 1, $false, "^p", 2)

# ---------------------------------------------------------------------
# 4. Re-create "_GoBack" on the new trailing empty paragraph, matching
#    where Word left the caret after the edit.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$pr = $lastPara.Range
$pr.Collapse(0)
$d.Bookmarks.Add("_GoBack", $pr) | Out-Null

Write-Host "Paragraphs:" $d.Paragraphs.Count
Write-Host $d.Content.Text
